# Weekly update: insert a new latest-week record at the top of the data
# block (row 38) for "Agrícola del Norte S.A. de Arica - Guayaba", pushing
# all existing rows (38-67) down by one (to 39-68).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 38; rows 38..67 shift down to 39..68.
$ws.Rows("38").Insert()

# Populate the newly inserted row 38 with the latest week's record.
$ws.Cells.Item(38, 1).Value  = 1
$ws.Cells.Item(38, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(38, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(38, 4).Value  = 45079
$ws.Cells.Item(38, 5).Value  = 15
$ws.Cells.Item(38, 6).Value  = "Fruta"
$ws.Cells.Item(38, 7).Value  = 100108
$ws.Cells.Item(38, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(38, 9).Value  = 100108001
$ws.Cells.Item(38, 10).Value = "Guayaba"
$ws.Cells.Item(38, 11).Value = "Sin especificar"
$ws.Cells.Item(38, 12).Value = "Primera"
$ws.Cells.Item(38, 13).Value = 200
$ws.Cells.Item(38, 14).Value = 5000
$ws.Cells.Item(38, 15).Value = 6000
$ws.Cells.Item(38, 16).Value = 5500
$ws.Cells.Item(38, 17).Value = "$/caja 10 kilos"
$ws.Cells.Item(38, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(38, 19).Value = 550
$ws.Cells.Item(38, 20).Value = 10
